$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Thailand/BBK row's Status columns (H2 and I2) flip from "Active" to "Inactive"
$ws.Range("H2").Value = "Inactive"
$ws.Range("I2").Value = "Inactive"

# Active cell selection moves from J2 to I3
[void]$ws.Range("I3").Select()
